# Update "want to go" counts (column F) for several rows on the
# "展览" (Exhibition) sheet and the mirrored "全部类型" (All Types) sheet,
# matching the refreshed data output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{Row = 9;  Value = 549},
    @{Row = 11; Value = 567},
    @{Row = 13; Value = 13463},
    @{Row = 14; Value = 179},
    @{Row = 17; Value = 5546},
    @{Row = 19; Value = 53}
)

$ws1 = $wb.Worksheets.Item("展览")
foreach ($u in $updates) {
    $ws1.Cells.Item($u.Row, 6).Value = $u.Value
}

$ws4 = $wb.Worksheets.Item("全部类型")
$updates4 = @(
    @{Row = 31; Value = 549},
    @{Row = 33; Value = 567},
    @{Row = 35; Value = 13463},
    @{Row = 36; Value = 179},
    @{Row = 40; Value = 5546},
    @{Row = 42; Value = 53}
)
foreach ($u in $updates4) {
    $ws4.Cells.Item($u.Row, 6).Value = $u.Value
}
